# Auto-generated edit script applying Rafflesia_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1111
$ws.Range("I40").Value = 1111
$ws.Range("K40").Value = 1111
$ws.Range("M40").Value = -936
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 935.125
$ws.Range("J80").Value = 925.8570999999999
$ws.Range("L80").Value = 2777.5713
$ws.Range("N80").Value = -4773.5713
$ws.Range("H83").Value = 935.125
$ws.Range("J83").Value = 925.8570999999999
$ws.Range("L83").Value = 8332.713899999999
$ws.Range("N83").Value = -18316.7139
$ws.Range("H106").Value = 2599
$ws.Range("I106").Value = 1998.75
$ws.Range("K106").Value = 1998.75
$ws.Range("M106").Value = -1367.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 15006
$ws.Range("J3").Value = 15006
$ws.Range("L3").Value = 15006
$ws.Range("N3").Value = -15236
$ws.Range("H32").Value = 7083.2666
$ws.Range("I32").Value = 6160.643
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 6160.643
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -5873.643
$ws.Range("N32").Value = -20574
$ws.Range("H122").Value = 518.3333
$ws.Range("I122").Value = 518.3333
$ws.Range("K122").Value = 1554.9999
$ws.Range("M122").Value = 895.0001
$ws.Range("H132").Value = 5500.3335
$ws.Range("I132").Value = 2750.5
$ws.Range("K132").Value = 8251.5
$ws.Range("M132").Value = -5721.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2542.25
$ws.Range("I20").Value = 2084.5
$ws.Range("K20").Value = 2084.5
$ws.Range("M20").Value = -1837.5
$ws.Range("H86").Value = 2802.7778
$ws.Range("I86").Value = 1287.5
$ws.Range("J86").Value = 5833.3335
$ws.Range("K86").Value = 1287.5
$ws.Range("L86").Value = 5833.3335
$ws.Range("M86").Value = -164.5
$ws.Range("N86").Value = -8079.3335
$ws.Range("H89").Value = 2802.7778
$ws.Range("I89").Value = 1287.5
$ws.Range("J89").Value = 5833.3335
$ws.Range("K89").Value = 6437.5
$ws.Range("L89").Value = 29166.6675
$ws.Range("M89").Value = -821.5
$ws.Range("N89").Value = -40398.6675
$ws.Range("H94").Value = 3827.5715
$ws.Range("I94").Value = 899.5
$ws.Range("J94").Value = 4998.8
$ws.Range("K94").Value = 899.5
$ws.Range("L94").Value = 4998.8
$ws.Range("M94").Value = -448.5
$ws.Range("N94").Value = -5900.8
$ws.Range("H105").Value = 2182
$ws.Range("I105").Value = 1683
$ws.Range("K105").Value = 1683
$ws.Range("M105").Value = 64

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 724.1
$ws.Range("I22").Value = 580.25
$ws.Range("J22").Value = 1299.5
$ws.Range("K22").Value = 580.25
$ws.Range("L22").Value = 1299.5
$ws.Range("M22").Value = -230.25
$ws.Range("N22").Value = -1999.5
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 4663.8
$ws.Range("I99").Value = 4663.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4663.8
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3165.8
$ws.Range("N99").ClearContents()
$ws.Range("H121").Value = 80000
$ws.Range("J121").Value = 80000
$ws.Range("L121").Value = 80000
$ws.Range("N121").Value = -82620
$ws.Range("H122").Value = 2996.5
$ws.Range("J122").Value = 2996.5
$ws.Range("L122").Value = 8989.5
$ws.Range("N122").Value = -13889.5
$ws.Range("H126").Value = 4663.8
$ws.Range("I126").Value = 4663.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13991.4
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11521.4
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 5999.6665
$ws.Range("I134").Value = 5999.6665
$ws.Range("K134").Value = 17998.9995
$ws.Range("M134").Value = -15463.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 7208.5
$ws.Range("I6").Value = 2408
$ws.Range("K6").Value = 2408
$ws.Range("M6").Value = -2295
$ws.Range("H16").Value = 7208.5
$ws.Range("I16").Value = 2408
$ws.Range("K16").Value = 2408
$ws.Range("M16").Value = -2158
$ws.Range("H32").Value = 19998
$ws.Range("I32").Value = 19996
$ws.Range("K32").Value = 19996
$ws.Range("M32").Value = -19700
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H70").Value = 7749.5
$ws.Range("I70").Value = 7000
$ws.Range("K70").Value = 7000
$ws.Range("M70").Value = -6730
$ws.Range("H73").Value = 7749.5
$ws.Range("I73").Value = 7000
$ws.Range("K73").Value = 7000
$ws.Range("M73").Value = -6064
$ws.Range("H80").Value = 23081.092
$ws.Range("I80").Value = 18511.25
$ws.Range("J80").Value = 25692.428
$ws.Range("K80").Value = 18511.25
$ws.Range("L80").Value = 25692.428
$ws.Range("M80").Value = -17513.25
$ws.Range("N80").Value = -27688.428
$ws.Range("H83").Value = 23081.092
$ws.Range("I83").Value = 18511.25
$ws.Range("J83").Value = 25692.428
$ws.Range("K83").Value = 92556.25
$ws.Range("L83").Value = 128462.14
$ws.Range("M83").Value = -87564.25
$ws.Range("N83").Value = -138446.14
$ws.Range("H97").Value = 1496
$ws.Range("I97").Value = 1232
$ws.Range("J97").Value = 2156
$ws.Range("K97").Value = 1232
$ws.Range("L97").Value = 2156
$ws.Range("M97").Value = -736
$ws.Range("N97").Value = -3148

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3013
$ws.Range("I16").Value = 4536.2
$ws.Range("J16").Value = 1489.8
$ws.Range("K16").Value = 4536.2
$ws.Range("L16").Value = 1489.8
$ws.Range("M16").Value = -4366.2
$ws.Range("N16").Value = -1829.8
$ws.Range("H40").Value = 999
$ws.Range("I40").Value = 999
$ws.Range("K40").Value = 999
$ws.Range("M40").Value = -863
$ws.Range("H82").Value = 600
$ws.Range("J82").Value = 600
$ws.Range("L82").Value = 600
$ws.Range("N82").Value = -1322
$ws.Range("H85").Value = 600
$ws.Range("J85").Value = 600
$ws.Range("L85").Value = 600
$ws.Range("N85").Value = -3096
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H122").Value = 8465
$ws.Range("I122").Value = 3400
$ws.Range("J122").Value = 10997.5
$ws.Range("K122").Value = 10200
$ws.Range("L122").Value = 32992.5
$ws.Range("M122").Value = -7750
$ws.Range("N122").Value = -37892.5
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 24499.5
$ws.Range("I61").Value = 24499.5
$ws.Range("K61").Value = 24499.5
$ws.Range("M61").Value = -24207.5

